# Update league results: recompute "Puntos" (points) for several rows and
# mark a batch of previously-timed competitors as disqualified
# ("Descalificado") across all four sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: Elite_Masc
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Elite_Masc")

# Recalculated points for rows 3-13 (times unchanged)
$ws.Cells.Item(3, 7).Value = 92
$ws.Cells.Item(4, 7).Value = 86
$ws.Cells.Item(5, 7).Value = 82
$ws.Cells.Item(6, 7).Value = 80
$ws.Cells.Item(7, 7).Value = 79
$ws.Cells.Item(8, 7).Value = 78
$ws.Cells.Item(9, 7).Value = 77
$ws.Cells.Item(10, 7).Value = 76
$ws.Cells.Item(11, 7).Value = 75
$ws.Cells.Item(12, 7).Value = 74
$ws.Cells.Item(13, 7).Value = 73

# Rows already shown as "00:00:00" stay disqualified, just re-labelled
$ws.Cells.Item(14, 6).Value = "Descalificado"
$ws.Cells.Item(15, 6).Value = "Descalificado"

# Rows 16-32: newly disqualified competitors (time -> "Descalificado", points -> 0)
$ws.Cells.Item(16, 6).Value = "Descalificado"
$ws.Cells.Item(16, 7).Value = 0
$ws.Cells.Item(17, 6).Value = "Descalificado"
$ws.Cells.Item(17, 7).Value = 0
$ws.Cells.Item(18, 6).Value = "Descalificado"
$ws.Cells.Item(18, 7).Value = 0
$ws.Cells.Item(19, 6).Value = "Descalificado"
$ws.Cells.Item(19, 7).Value = 0
$ws.Cells.Item(20, 6).Value = "Descalificado"
$ws.Cells.Item(20, 7).Value = 0
$ws.Cells.Item(21, 6).Value = "Descalificado"
$ws.Cells.Item(21, 7).Value = 0
$ws.Cells.Item(22, 6).Value = "Descalificado"
$ws.Cells.Item(22, 7).Value = 0
$ws.Cells.Item(23, 6).Value = "Descalificado"
$ws.Cells.Item(23, 7).Value = 0
$ws.Cells.Item(24, 6).Value = "Descalificado"
$ws.Cells.Item(24, 7).Value = 0
$ws.Cells.Item(25, 6).Value = "Descalificado"
$ws.Cells.Item(25, 7).Value = 0
$ws.Cells.Item(26, 6).Value = "Descalificado"
$ws.Cells.Item(26, 7).Value = 0
$ws.Cells.Item(27, 6).Value = "Descalificado"
$ws.Cells.Item(27, 7).Value = 0
$ws.Cells.Item(28, 6).Value = "Descalificado"
$ws.Cells.Item(28, 7).Value = 0
$ws.Cells.Item(29, 6).Value = "Descalificado"
$ws.Cells.Item(29, 7).Value = 0
$ws.Cells.Item(30, 6).Value = "Descalificado"
$ws.Cells.Item(30, 7).Value = 0
$ws.Cells.Item(31, 6).Value = "Descalificado"
$ws.Cells.Item(31, 7).Value = 0
$ws.Cells.Item(32, 6).Value = "Descalificado"
$ws.Cells.Item(32, 7).Value = 0

# ---------------------------------------------------------------------------
# Sheet: Elite_Fem
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Elite_Fem")

# Recalculated points for rows 2-4 (times unchanged)
$ws.Cells.Item(2, 7).Value = 100
$ws.Cells.Item(3, 7).Value = 92
$ws.Cells.Item(4, 7).Value = 86

# Rows 5-11: newly disqualified competitors (time -> "Descalificado", points -> 0)
$ws.Cells.Item(5, 6).Value = "Descalificado"
$ws.Cells.Item(5, 7).Value = 0
$ws.Cells.Item(6, 6).Value = "Descalificado"
$ws.Cells.Item(6, 7).Value = 0
$ws.Cells.Item(7, 6).Value = "Descalificado"
$ws.Cells.Item(7, 7).Value = 0
$ws.Cells.Item(8, 6).Value = "Descalificado"
$ws.Cells.Item(8, 7).Value = 0
$ws.Cells.Item(9, 6).Value = "Descalificado"
$ws.Cells.Item(9, 7).Value = 0
$ws.Cells.Item(10, 6).Value = "Descalificado"
$ws.Cells.Item(10, 7).Value = 0
$ws.Cells.Item(11, 6).Value = "Descalificado"
$ws.Cells.Item(11, 7).Value = 0

# ---------------------------------------------------------------------------
# Sheet: GGEE_Masc
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("GGEE_Masc")

# Recalculated points for rows 2-9 (times unchanged)
$ws.Cells.Item(2, 7).Value = 100
$ws.Cells.Item(3, 7).Value = 98.09999999999999
$ws.Cells.Item(4, 7).Value = 90.59999999999999
$ws.Cells.Item(5, 7).Value = 90.5
$ws.Cells.Item(6, 7).Value = 86.5
$ws.Cells.Item(7, 7).Value = 85.90000000000001
$ws.Cells.Item(8, 7).Value = 80.8
$ws.Cells.Item(9, 7).Value = 64.40000000000001

# Rows already shown as "00:00:00" stay disqualified, just re-labelled
$ws.Cells.Item(10, 6).Value = "Descalificado"
$ws.Cells.Item(11, 6).Value = "Descalificado"

# Rows 12-30: newly disqualified competitors (time -> "Descalificado", points -> 0)
$ws.Cells.Item(12, 6).Value = "Descalificado"
$ws.Cells.Item(12, 7).Value = 0
$ws.Cells.Item(13, 6).Value = "Descalificado"
$ws.Cells.Item(13, 7).Value = 0
$ws.Cells.Item(14, 6).Value = "Descalificado"
$ws.Cells.Item(14, 7).Value = 0
$ws.Cells.Item(15, 6).Value = "Descalificado"
$ws.Cells.Item(15, 7).Value = 0
$ws.Cells.Item(16, 6).Value = "Descalificado"
$ws.Cells.Item(16, 7).Value = 0
$ws.Cells.Item(17, 6).Value = "Descalificado"
$ws.Cells.Item(17, 7).Value = 0
$ws.Cells.Item(18, 6).Value = "Descalificado"
$ws.Cells.Item(18, 7).Value = 0
$ws.Cells.Item(19, 6).Value = "Descalificado"
$ws.Cells.Item(19, 7).Value = 0
$ws.Cells.Item(20, 6).Value = "Descalificado"
$ws.Cells.Item(20, 7).Value = 0
$ws.Cells.Item(21, 6).Value = "Descalificado"
$ws.Cells.Item(21, 7).Value = 0
$ws.Cells.Item(22, 6).Value = "Descalificado"
$ws.Cells.Item(22, 7).Value = 0
$ws.Cells.Item(23, 6).Value = "Descalificado"
$ws.Cells.Item(23, 7).Value = 0
$ws.Cells.Item(24, 6).Value = "Descalificado"
$ws.Cells.Item(24, 7).Value = 0
$ws.Cells.Item(25, 6).Value = "Descalificado"
$ws.Cells.Item(25, 7).Value = 0
$ws.Cells.Item(26, 6).Value = "Descalificado"
$ws.Cells.Item(26, 7).Value = 0
$ws.Cells.Item(27, 6).Value = "Descalificado"
$ws.Cells.Item(27, 7).Value = 0
$ws.Cells.Item(28, 6).Value = "Descalificado"
$ws.Cells.Item(28, 7).Value = 0
$ws.Cells.Item(29, 6).Value = "Descalificado"
$ws.Cells.Item(29, 7).Value = 0
$ws.Cells.Item(30, 6).Value = "Descalificado"
$ws.Cells.Item(30, 7).Value = 0

# ---------------------------------------------------------------------------
# Sheet: GGEE_Fem
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("GGEE_Fem")

# Rows 2-4: newly disqualified competitors (time -> "Descalificado", points -> 0)
$ws.Cells.Item(2, 6).Value = "Descalificado"
$ws.Cells.Item(2, 7).Value = 0
$ws.Cells.Item(3, 6).Value = "Descalificado"
$ws.Cells.Item(3, 7).Value = 0
$ws.Cells.Item(4, 6).Value = "Descalificado"
$ws.Cells.Item(4, 7).Value = 0
